$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(427.9782062899712,3374.105,-803.177),
    @(456.3503240472388,3010.388,-681.174),
    @(376.2774675207996,2937.464,-679.5700000000001),
    @(408.5741938640421,3205.226,-746.95),
    @(387.0446386159527,3921.591,-907.098),
    @(369.1798734478832,3019.653,-668.325),
    @(312.9227364959329,3485.01,-800.76),
    @(392.2561362963963,3680.075,-864.307),
    @(380.7872620768668,3995.136,-939.732),
    @(530.8715230473686,3411.079,-779.297),
    @(382.3965490247278,3053.366,-722.427),
    @(410.730922019671,3290.369,-736.924),
    @(533.4681356085605,3302.031,-753.702),
    @(430.3367894376382,3472.779,-818.35),
    @(674.6629003601136,3874.8,-897.389),
    @(438.0092725114826,3531.565,-814.419),
    @(688.4379269308076,3628.003,-802.16),
    @(492.2175755131212,3570.881,-815.899),
    @(445.9809769653837,4237.553,-989.367),
    @(482.9289472972656,3801.93,-900.956),
    @(325.9809711756966,3125.726,-722.312),
    @(374.3476748999545,5390.704,-1291.607),
    @(446.8761000313985,2962.411,-685.606),
    @(621.4420870615282,3835.514,-901.747),
    @(388.2866908372181,3850.456,-875.9160000000001),
    @(307.355801056705,3902.123,-888.718),
    @(471.1023505138583,3606.27,-849.633),
    @(398.1929319938309,2873.417,-665.151),
    @(366.8861398576823,3050.101,-701.407),
    @(421.8554205074773,4139.028,-978.626),
    @(683.8419815910842,6702.29,-1612.015),
    @(493.5397201860787,3067.66,-695.314),
    @(350.0471442470368,3681.09,-794.549),
    @(474.8455402235411,3948.737,-920.131),
    @(581.6796592647529,3066.572,-719.04),
    @(484.5252051276414,3374.393,-766.361),
    @(377.5864124150023,3915.663,-909.831),
    @(555.0158890242848,3062.401,-699.165),
    @(384.8319075359489,3540.577,-803.572),
    @(439.6538973941318,4434.427,-1050.022),
    @(338.8121586826168,3366.851,-748.418),
    @(392.3103384902411,4163.534,-959.072),
    @(395.1749857726152,3126.082,-727.221),
    @(394.8568349133555,3594.55,-819.454),
    @(433.6155259941552,3020.47,-692.724),
    @(411.8166612227795,3366.948,-790.314),
    @(424.0009022697134,4254.943,-993.772),
    @(454.7426359238212,3511.724,-816.1559999999999),
    @(415.7674951144306,3376.009,-792.9589999999999),
    @(424.010216246149,3185.493,-727.229)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

Write-Output "Updated A2:C51 with new values"